# Read from email logic added
# ------------------------------------------------------------------
# 1. Update row 2 of the "data" sheet with the new comment content
#    (subject/body/sender changed; sender hyperlink re-pointed).
# 2. Add a new worksheet "Sheet1" after "data" containing the new
#    "Home Insurance Claim" rows (rows 2 and 8) plus the same updated
#    comment as row 5 of the "data" sheet (with matching formatting).
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("data")

# --- Update data!A2:H2 -------------------------------------------------
$ws1.Range("A2").Value = "A1"
$ws1.Range("B2").Value = "Expected better service again"
$ws1.Range("C2").Value = "Room was a bit dark, facilities could have been much better. Will not visit again."
$ws1.Range("D2").Value = "leslei201@mitchell.com"
$ws1.Range("E2").Value = "berniece.heller@schamberger.org"
$ws1.Range("F2").Value = "vilma32@monahan.com"
$ws1.Range("G2").Value = "2023-06-19T10:33:12Z"
$ws1.Range("H2").Value = "Bronze"

# Re-point the sender hyperlink to the new e-mail address.
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("D2"), "mailto:leslei201@mitchell.com")
$ws1.Range("D2").Style = "Hyperlink"

# Row got shorter (single-line subject/body instead of the long review).
$ws1.Rows(2).RowHeight = 29

# --- Add the new "Sheet1" worksheet, placed after "data" ---------------
$ws2 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws2.Name = "Sheet1"

# Row 2 - Home Insurance Claim (amber@customer.com -> underwriter@insurer.com)
$ws2.Range("A2").Value = "A1"
$ws2.Range("B2").Value = "Home Insurance Claim"
$ws2.Range("C2").Value = "I would like to file a claim for earthquake damage to my home."
$ws2.Range("D2").Value = "amber@customer.com"
$ws2.Range("E2").Value = "underwriter@insurer.com"
$ws2.Range("F2").Value = "vilma32@monahan.com"
$ws2.Range("G2").Value = "2023-06-19T10:33:12Z"
$ws2.Range("H2").Value = "Bronze"

# Row 5 - same content/formatting as data!A2:H2
$ws2.Range("A5").Value = "A1"
$ws2.Range("B5").Value = "Expected better service again"
$ws2.Range("C5").Value = "Room was a bit dark, facilities could have been much better. Will not visit again."
$ws2.Range("D5").Value = "leslei201@mitchell.com"
$ws2.Range("E5").Value = "berniece.heller@schamberger.org"
$ws2.Range("F5").Value = "vilma32@monahan.com"
$ws2.Range("G5").Value = "2023-06-19T10:33:12Z"
$ws2.Range("H5").Value = "Bronze"

$ws2.Range("B5:C5").WrapText = $true
$ws2.Hyperlinks.Add($ws2.Range("D5"), "mailto:leslei201@mitchell.com")
$ws2.Range("D5").Style = "Hyperlink"
$ws2.Rows(5).RowHeight = 72

# Row 8 - Home Insurance Claim again
$ws2.Range("A8").Value = "A1"
$ws2.Range("B8").Value = "Home Insurance Claim"
$ws2.Range("C8").Value = "I would like to file a claim for earthquake damage to my home."
$ws2.Range("D8").Value = "amber@customer.com"
$ws2.Range("E8").Value = "underwriter@insurer.com"
$ws2.Range("F8").Value = "vilma32@monahan.com"
$ws2.Range("G8").Value = "2023-06-19T10:33:12Z"
$ws2.Range("H8").Value = "Bronze"

# Match the saved selection on Sheet1 (A5:H5 selected, H5 active).
$ws2.Select() | Out-Null
$ws2.Range("A5:H5").Activate()

# Leave "data" as the active sheet, as it was before the edit.
$ws1.Activate()
$ws1.Range("A2").Select() | Out-Null
